# Demo Verification Script fixes
# Updates the "DateProd" (column B) run-timestamp cells on the various
# VT-*-Generic worksheets to the latest verification-run timestamps.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "VT-P-DebitVoid-DualCF-Generic";    Cell = "B2"; Value = "Tue Aug 12 02:23:36 IST 2025" },

    @{ Sheet = "VT-P-DebitVoid-SingleCF-Generic";  Cell = "B2"; Value = "Tue Aug 12 02:28:22 IST 2025" },
    @{ Sheet = "VT-P-DebitVoid-SingleCF-Generic";  Cell = "B3"; Value = "Tue Aug 12 02:29:18 IST 2025" },
    @{ Sheet = "VT-P-DebitVoid-SingleCF-Generic";  Cell = "B4"; Value = "Tue Aug 12 02:30:15 IST 2025" },
    @{ Sheet = "VT-P-DebitVoid-SingleCF-Generic";  Cell = "B5"; Value = "Tue Aug 12 02:31:14 IST 2025" },

    @{ Sheet = "VT-P-DebitVoid-NoCF-Generic";      Cell = "B2"; Value = "Tue Aug 12 02:24:33 IST 2025" },
    @{ Sheet = "VT-P-DebitVoid-NoCF-Generic";      Cell = "B3"; Value = "Tue Aug 12 02:25:29 IST 2025" },
    @{ Sheet = "VT-P-DebitVoid-NoCF-Generic";      Cell = "B4"; Value = "Tue Aug 12 02:26:26 IST 2025" },
    @{ Sheet = "VT-P-DebitVoid-NoCF-Generic";      Cell = "B5"; Value = "Tue Aug 12 02:27:19 IST 2025" },

    @{ Sheet = "VT-P-DebitCredit-DualCF-Generic";  Cell = "B2"; Value = "Tue Aug 12 02:12:19 IST 2025" },
    @{ Sheet = "VT-P-DebitCredit-DualCF-Generic";  Cell = "B3"; Value = "Tue Aug 12 02:13:15 IST 2025" },
    @{ Sheet = "VT-P-DebitCredit-DualCF-Generic";  Cell = "B4"; Value = "Tue Aug 12 02:14:11 IST 2025" },
    @{ Sheet = "VT-P-DebitCredit-DualCF-Generic";  Cell = "B5"; Value = "Tue Aug 12 02:15:14 IST 2025" },

    @{ Sheet = "VT-P-DebitCredit-SingleCF-Gener";  Cell = "B2"; Value = "Tue Aug 12 02:19:46 IST 2025" },
    @{ Sheet = "VT-P-DebitCredit-SingleCF-Gener";  Cell = "B3"; Value = "Tue Aug 12 02:20:43 IST 2025" },
    @{ Sheet = "VT-P-DebitCredit-SingleCF-Gener";  Cell = "B4"; Value = "Tue Aug 12 02:21:39 IST 2025" },
    @{ Sheet = "VT-P-DebitCredit-SingleCF-Gener";  Cell = "B5"; Value = "Tue Aug 12 02:22:36 IST 2025" },

    @{ Sheet = "VT-P-DebitCredit-NoCF-Generic";    Cell = "B2"; Value = "Tue Aug 12 02:16:08 IST 2025" },
    @{ Sheet = "VT-P-DebitCredit-NoCF-Generic";    Cell = "B3"; Value = "Tue Aug 12 02:17:02 IST 2025" },
    @{ Sheet = "VT-P-DebitCredit-NoCF-Generic";    Cell = "B4"; Value = "Tue Aug 12 02:17:56 IST 2025" },
    @{ Sheet = "VT-P-DebitCredit-NoCF-Generic";    Cell = "B5"; Value = "Tue Aug 12 02:18:53 IST 2025" },

    @{ Sheet = "VT-C-DebitCredit-DualCF-Generic";  Cell = "B2"; Value = "Tue Aug 12 02:09:33 IST 2025" },
    @{ Sheet = "VT-C-DebitCredit-SingleCF-Gener";  Cell = "B2"; Value = "Tue Aug 12 02:11:23 IST 2025" },
    @{ Sheet = "VT-C-DebitCredit-NoCF-Generic";    Cell = "B2"; Value = "Tue Aug 12 02:10:31 IST 2025" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
